$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.65897045142951
$ws.Range("C2").Value = 9.459933672100128
$ws.Range("D2").Value = 5.963794528825298
$ws.Range("E2").Value = 11.66460203902936
$ws.Range("G2").Value = 22.07014778229567
$ws.Range("H2").Value = 12.42001788436983
$ws.Range("M2").Value = 14.42093827188041
$ws.Range("O2").Value = 18.03175421150882
$ws.Range("B3").Value = 11.98859368823048
$ws.Range("C3").Value = 9.045193081341649
$ws.Range("D3").Value = 5.84102843233049
$ws.Range("E3").Value = 11.59755345091735
$ws.Range("G3").Value = 22.03665919583081
$ws.Range("H3").Value = 12.47531389426312
$ws.Range("M3").Value = 14.09635239970225
$ws.Range("O3").Value = 18.10323093379353
$ws.Range("B4").Value = 11.55731064408175
$ws.Range("C4").Value = 8.779511201313406
$ws.Range("D4").Value = 5.766121091900955
$ws.Range("E4").Value = 11.56152802797377
$ws.Range("G4").Value = 22.02904870656905
$ws.Range("H4").Value = 12.51233423791978
$ws.Range("M4").Value = 13.89591426323856
$ws.Range("O4").Value = 18.15356921308154
$ws.Range("B5").Value = 11.37677564725525
$ws.Range("C5").Value = 8.668577583644046
$ws.Range("D5").Value = 5.735762273491301
$ws.Range("E5").Value = 11.5481510330676
$ws.Range("G5").Value = 22.02919405111814
$ws.Range("H5").Value = 12.52818968106621
$ws.Range("M5").Value = 13.81406445536871
$ws.Range("O5").Value = 18.17569478256113
$ws.Range("B6").Value = 11.34651438481357
$ws.Range("C6").Value = 8.649999570897684
$ws.Range("D6").Value = 5.7307326859775
$ws.Range("E6").Value = 11.54600878152341
$ws.Range("G6").Value = 22.02941386790678
$ws.Range("H6").Value = 12.53086886734243
$ws.Range("M6").Value = 13.80046661312211
$ws.Range("O6").Value = 18.17946579829338
$ws.Range("B7").Value = 11.55489501764509
$ws.Range("C7").Value = 8.778025752975307
$ws.Range("D7").Value = 5.765710924154926
$ws.Range("E7").Value = 11.56134233095955
$ws.Range("G7").Value = 22.02903753877694
$ws.Range("H7").Value = 12.51254495785092
$ws.Range("M7").Value = 13.8948109371973
$ws.Range("O7").Value = 18.15386109151022
$ws.Range("B8").Value = 12.43199498253229
$ws.Range("C8").Value = 9.319279838392799
$ws.Range("D8").Value = 5.921395463431798
$ws.Range("E8").Value = 11.64042348373847
$ws.Range("G8").Value = 22.05590851761061
$ws.Range("H8").Value = 12.43844578265318
$ws.Range("M8").Value = 14.30932703093185
$ws.Range("O8").Value = 18.05505458538661
$ws.Range("B9").Value = 13.99021747713159
$ws.Range("C9").Value = 10.28921307235103
$ws.Range("D9").Value = 6.228370773790819
$ws.Range("E9").Value = 11.83570782789968
$ws.Range("G9").Value = 22.21152713789743
$ws.Range("H9").Value = 12.31758200107652
$ws.Range("M9").Value = 15.1082028439574
$ws.Range("O9").Value = 17.9129068449285
$ws.Range("B10").Value = 15.03044014806999
$ws.Range("C10").Value = 10.94171100743411
$ws.Range("D10").Value = 6.452260049337529
$ws.Range("E10").Value = 12.00270473462166
$ws.Range("G10").Value = 22.38840043135468
$ws.Range("H10").Value = 12.24381892925969
$ws.Range("M10").Value = 15.68036201940441
$ws.Range("O10").Value = 17.84048266990012
$ws.Range("B11").Value = 15.48007005566707
$ws.Range("C11").Value = 11.22480573725516
$ws.Range("D11").Value = 6.553226516308787
$ws.Range("E11").Value = 12.08352098866325
$ws.Range("G11").Value = 22.4822725066457
$ws.Range("H11").Value = 12.21355533614051
$ws.Range("M11").Value = 15.93623843238516
$ws.Range("O11").Value = 17.81459238952311
$ws.Range("B12").Value = 15.64688920141277
$ws.Range("C12").Value = 11.32998669261301
$ws.Range("D12").Value = 6.591291961171597
$ws.Range("E12").Value = 12.11479598950899
$ws.Range("G12").Value = 22.51972453922826
$ws.Range("H12").Value = 12.20257096069177
$ws.Range("M12").Value = 16.03240735706451
$ws.Range("O12").Value = 17.80581083789381
$ws.Range("B13").Value = 15.61111574152661
$ws.Range("C13").Value = 11.30742462765512
$ws.Range("D13").Value = 6.583102010110968
$ws.Range("E13").Value = 12.10803091398096
$ws.Range("G13").Value = 22.51157430164469
$ws.Range("H13").Value = 12.20491544267123
$ws.Range("M13").Value = 16.01172935353582
$ws.Range("O13").Value = 17.80765650555152
$ws.Range("B14").Value = 15.49386370852205
$ws.Range("C14").Value = 11.23349977736909
$ws.Range("D14").Value = 6.55636178693659
$ws.Range("E14").Value = 12.08608069422606
$ws.Range("G14").Value = 22.4853156915105
$ws.Range("H14").Value = 12.21264209541318
$ws.Range("M14").Value = 15.94416536509789
$ws.Range("O14").Value = 17.8138493898405
$ws.Range("B15").Value = 15.42159322624299
$ws.Range("C15").Value = 11.18795424987806
$ws.Range("D15").Value = 6.539959510626549
$ws.Range("E15").Value = 12.07272222726884
$ws.Range("G15").Value = 22.46947877911505
$ws.Range("H15").Value = 12.21743692786387
$ws.Range("M15").Value = 15.9026832893648
$ws.Range("O15").Value = 17.8177760983655
$ws.Range("B16").Value = 15.00057634726991
$ws.Range("C16").Value = 10.92292965717685
$ws.Range("D16").Value = 6.445640300160739
$ws.Range("E16").Value = 11.99751862450964
$ws.Range("G16").Value = 22.38253363802509
$ws.Range("H16").Value = 12.24586317278744
$ws.Range("M16").Value = 15.66354322300925
$ws.Range("O16").Value = 17.84231733346991
$ws.Range("B17").Value = 14.73621463087675
$ws.Range("C17").Value = 10.75679337425447
$ws.Range("D17").Value = 6.387522318171474
$ws.Range("E17").Value = 11.95260789210487
$ws.Range("G17").Value = 22.33261549930287
$ws.Range("H17").Value = 12.26414657185286
$ws.Range("M17").Value = 15.5156429201584
$ws.Range("O17").Value = 17.85918559144466
$ws.Range("B18").Value = 14.5819444969767
$ws.Range("C18").Value = 10.65994608242504
$ws.Range("D18").Value = 6.354013720701418
$ws.Range("E18").Value = 11.92723445311546
$ws.Range("G18").Value = 22.3051680657807
$ws.Range("H18").Value = 12.27497237457689
$ws.Range("M18").Value = 15.43016541379988
$ws.Range("O18").Value = 17.86955146262234
$ws.Range("B19").Value = 14.52933255121876
$ws.Range("C19").Value = 10.62693523056996
$ws.Range("D19").Value = 6.342655735644233
$ws.Range("E19").Value = 11.91872285687388
$ws.Range("G19").Value = 22.29609258882325
$ws.Range("H19").Value = 12.27869091689821
$ws.Range("M19").Value = 15.40115686572295
$ws.Range("O19").Value = 17.87317490715177
$ws.Range("B20").Value = 14.76458620808914
$ws.Range("C20").Value = 10.7746127480129
$ws.Range("D20").Value = 6.393717734477324
$ws.Range("E20").Value = 11.95734149144466
$ws.Range("G20").Value = 22.33779868900229
$ws.Range("H20").Value = 12.26216820463889
$ws.Range("M20").Value = 15.53143020128316
$ws.Range("O20").Value = 17.85732118718378
$ws.Range("B21").Value = 15.52839738347634
$ws.Range("C21").Value = 11.25526847524266
$ws.Range("D21").Value = 6.564220927040799
$ws.Range("E21").Value = 12.09250999075244
$ws.Range("G21").Value = 22.49297700681524
$ws.Range("H21").Value = 12.21035965790939
$ws.Range("M21").Value = 15.96403095549097
$ws.Range("O21").Value = 17.81200257844744
$ws.Range("B22").Value = 16.00748851609906
$ws.Range("C22").Value = 11.55760938592629
$ws.Range("D22").Value = 6.674655738127812
$ws.Range("E22").Value = 12.18475212994594
$ws.Range("G22").Value = 22.60548210570309
$ws.Range("H22").Value = 12.17927416941034
$ws.Range("M22").Value = 16.24249090878688
$ws.Range("O22").Value = 17.78834757069653
$ws.Range("B23").Value = 15.75364361445155
$ws.Range("C23").Value = 11.39733699576193
$ws.Range("D23").Value = 6.615819087810118
$ws.Range("E23").Value = 12.13517267309291
$ws.Range("G23").Value = 22.54443088951116
$ws.Range("H23").Value = 12.19561041688606
$ws.Range("M23").Value = 16.09429126008405
$ws.Range("O23").Value = 17.80042460502352
$ws.Range("B24").Value = 14.75176653112965
$ws.Range("C24").Value = 10.76656075860911
$ws.Range("D24").Value = 6.390917081981621
$ws.Range("E24").Value = 11.95520003895659
$ws.Range("G24").Value = 22.33545146985541
$ws.Range("H24").Value = 12.26306164516667
$ws.Range("M24").Value = 15.5242941586884
$ws.Range("O24").Value = 17.85816200386405
$ws.Range("B25").Value = 13.5866929937429
$ws.Range("C25").Value = 10.03710845139549
$ws.Range("D25").Value = 6.14542819341815
$ws.Range("E25").Value = 11.77866595112076
$ws.Range("G25").Value = 22.15840653663512
$ws.Range("H25").Value = 12.347648048132
$ws.Range("M25").Value = 14.89425040265501
$ws.Range("O25").Value = 17.94577759096846
